$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify D1's rich-text value down to plain "LCSC Part #"
$ws.Range("D1").Value = "LCSC Part #"

# Update the active selection from A9 to A17
$ws.Range("A17").Select()
